# "Updated general results excel"
#
# The original sheet held a single 8-row results table (header in row 1,
# six dimension rows, a "Total" row). This edit:
#   1) inserts a new row above everything and labels it "ar" (so the old
#      table now lives in rows 2-9 instead of 1-8);
#   2) adds a second, parallel results block further down (rows 12-18) that
#      mirrors the header/dimension layout but computes its values from a
#      (currently empty) raw-score area below, so every formula resolves to
#      #DIV/0! for now;
#   3) nudges the small helper text box on the sheet down by one row to
#      keep it aligned with the table it was sitting next to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row 1; everything else shifts down by one row. -------
$ws.Rows.Item(1).Insert() | Out-Null
$ws.Range("A1").Value = "ar"

# --- 2) Second block: header (row 13) + 5 dimension rows (14-18). --------
$ws.Range("A12").Value = "vr"

$ws.Range("B13").Value = "Utilizador 1"
$ws.Range("C13").Value = "Utilizador 2"
$ws.Range("D13").Value = "Utilizador 3"
$ws.Range("E13").Value = "Mean"

$ws.Range("A14").Value = "Total CF"
$ws.Range("A15").Value = "Total SF"
$ws.Range("A16").Value = "Total DF"
$ws.Range("A17").Value = "Total RF"
$ws.Range("A18").Value = "Total P"

# Row 14 (CF): references rows 21/27 individually (no fill-down pattern).
$ws.Range("B14").Formula = "=(B21-B27)/(C27-B27)"
$ws.Range("C14").Formula = "=(C21-B27)/(C27-B27)"
$ws.Range("D14").Formula = "=(D21-B27)/(C27-B27)"
$ws.Range("E14").Formula = "=MEDIAN(B14:D14)"

# Rows 15-18 (SF/DF/RF/P): each B/C cell is its own formula, while D and E
# are filled across the whole block in one shot so they come out as a
# shared formula group, same shape as the original E-column median block.
$ws.Range("B15").Formula = "=(B22-B28)/(C28-B28)"
$ws.Range("C15").Formula = "=(C22-B28)/(C28-B28)"

$ws.Range("B16").Formula = "=(B23-B29)/(C29-B29)"
$ws.Range("C16").Formula = "=(C23-B29)/(C29-B29)"

$ws.Range("B17").Formula = "=(B24-B30)/(C30-B30)"
$ws.Range("C17").Formula = "=(C24-B30)/(C30-B30)"

$ws.Range("B18").Formula = "=(B25-B31)/(C31-B31)"
$ws.Range("C18").Formula = "=(C25-B31)/(C31-B31)"

$ws.Range("D15:D18").Formula = "=(D22-B28)/(C28-B28)"
$ws.Range("E15:E18").Formula = "=MEDIAN(B15:D15)"

# Column E has no column-level style, so E13/E14 need the 2-decimal number
# format applied explicitly (same "0.00" style already used by columns B-D).
$ws.Range("E13:E18").NumberFormat = "0.00"

# --- 3) Move the little text-box shape down by one row (~16pt). ----------
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top + 16

# --- Restore the selection to match what was saved with the table. -------
$ws.Range("A2:E9").Select() | Out-Null
